# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" detail table (rows 16-42, columns C:F on Hoja1) is
# reshuffled: the existing worker (CC 45517938 - DERLY ZARATE LLERENA) and a
# newly-added worker (CC 71729664 - GABRIEL JAIME PAREJA) are interleaved
# row-by-row, each one listing all of their overdue periods in chronological
# order (1912, 2001..2012, 2101), instead of being grouped one-worker-per-block
# with the periods in descending order. The "Valor Mora" column (F) follows
# the period: partial periods 1912/2101 -> 19875/26500, full periods -> 33125.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

function Set-MoraRow($row, $doc, $nombre, $periodo, $valor) {
    $ws.Range("C$row").Value = $doc
    $ws.Range("D$row").Value = $nombre
    $ws.Range("E$row").Value = $periodo
    $ws.Range("F$row").Value = $valor
}

$derly   = "45517938"
$derlyNom = "DERLY ZARATE LLERENA"
$gabriel = "71729664"
$gabrielNom = "GABRIEL JAIME PAREJA"

Set-MoraRow 16 $gabriel $gabrielNom "1912" 19875

Set-MoraRow 17 $derly   $derlyNom   "2001" 33125
Set-MoraRow 18 $gabriel $gabrielNom "2001" 33125

Set-MoraRow 19 $derly   $derlyNom   "2002" 33125
Set-MoraRow 20 $gabriel $gabrielNom "2002" 33125

Set-MoraRow 21 $derly   $derlyNom   "2003" 33125
Set-MoraRow 22 $gabriel $gabrielNom "2003" 33125

Set-MoraRow 23 $derly   $derlyNom   "2004" 33125
Set-MoraRow 24 $gabriel $gabrielNom "2004" 33125

Set-MoraRow 25 $derly   $derlyNom   "2005" 33125
Set-MoraRow 26 $gabriel $gabrielNom "2005" 33125

Set-MoraRow 27 $derly   $derlyNom   "2006" 33125
Set-MoraRow 28 $gabriel $gabrielNom "2006" 33125

Set-MoraRow 29 $derly   $derlyNom   "2007" 33125
Set-MoraRow 30 $gabriel $gabrielNom "2007" 33125

Set-MoraRow 31 $derly   $derlyNom   "2008" 33125
Set-MoraRow 32 $gabriel $gabrielNom "2008" 33125

Set-MoraRow 33 $derly   $derlyNom   "2009" 33125
Set-MoraRow 34 $gabriel $gabrielNom "2009" 33125

Set-MoraRow 35 $derly   $derlyNom   "2010" 33125
Set-MoraRow 36 $gabriel $gabrielNom "2010" 33125

Set-MoraRow 37 $derly   $derlyNom   "2011" 33125
Set-MoraRow 38 $gabriel $gabrielNom "2011" 33125

Set-MoraRow 39 $derly   $derlyNom   "2012" 33125
Set-MoraRow 40 $gabriel $gabrielNom "2012" 33125

Set-MoraRow 41 $derly   $derlyNom   "2101" 26500
Set-MoraRow 42 $gabriel $gabrielNom "2101" 26500
